# Add new columns I ("I0") and J ("IF") to the active worksheet,
# filling in header labels and per-row numeric values for rows 2..62.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for the two new columns, matching the style used by the
# other header cells in row 1 (bold / bordered / centered "s=1" style).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Values for column I (I0) and column J (IF), rows 2 through 62.
$iVals = @(7,9,5,7,9,9,9,8,9,9,8,9,8,8,9,9,8,8,8,9,8,7,8,8,8,8,9,8,8,9,7,8,7,8,7,7,7,8,7,9,9,8,9,7,7,9,8,7,9,10,9,9,7,7,8,6,8,6,9,7,6)
$jVals = @(8,9,7,7,9,9,9,8,9,9,8,9,8,8,9,9,8,8,8,9,8,7,8,8,8,8,9,8,8,9,7,8,7,8,7,8,7,8,7,9,9,8,9,7,8,9,8,8,9,10,9,9,7,7,8,6,8,7,9,7,6)

for ($r = 2; $r -le 62; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 9).Value = $iVals[$idx]
    $ws.Cells.Item($r, 10).Value = $jVals[$idx]
}
